$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The table already has a blank/unused row sitting right after the "3 -
# put dev of account mgmt page into backlog" row. Turn that spare row into
# the new "4." change-log entry instead of inserting a brand-new row, so
# the overall row/column structure of the table is preserved.
$newRow = $t.Rows.Item(5)

$newRow.Cells.Item(1).Range.Text = "4. "

$newRow.Cells.Item(2).Range.Text = "Add user name to the user table in the database, the UserPersistance class and the User class"

$newRow.Cells.Item(3).Range.Text = "Åsa Wegelius"

$newRow.Cells.Item(4).Range.Text = "09-05-16"

$newRow.Cells.Item(5).Range.Text = "Approved"

$newRow.Cells.Item(6).Range.Text = "09-05-16"

$newRow.Cells.Item(7).Range.Text = "Users will be uncomfortable with having their email address in plain sight when they have logged in. Many prefer a user name."

Write-Output "done"
